$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for each data row (2-115).
# Bump every value in that range by one day: 46081 -> 46082.
$range = $ws.Range("C2:C115")
for ($i = 1; $i -le $range.Rows.Count; $i++) {
    $cell = $range.Cells.Item($i, 1)
    $cell.Value2 = $cell.Value2 + 1
}
